# Add total price calculator, stable
# - Rename header "product name " (B1) to "product"
# - Add row 13: HandWash product entry (A13 = Google Form URL, B13 = HandWash, C13 = 75)
# - Add row 14: Pencil product entry (A14 = VIT help center URL, B14 = Pencil, C14 = 25)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header in B1 from "product name " to "product"
$ws.Range("B1").Value = "product"

# New row 13 - HandWash (written in column order A, B, C to keep the
# shared-strings table ordered the same way Excel would produce it)
$ws.Range("A13").Value = "https://docs.google.com/forms/d/e/1FAIpQLSdiOsWQMywLXJV4NX9jYJwzXsW3mjDjmJDOE5EAsFzcN3N30g/viewform?usp=sf_link"
$ws.Range("B13").Value = "HandWash"
$ws.Range("C13").Value = 75

# New row 14 - Pencil (B14 written before A14 so "Pencil" precedes the
# help-center URL in the shared strings table)
$ws.Range("B14").Value = "Pencil"
$ws.Range("A14").Value = "https://vithelpcenter.vit.ac.in/vitcc-help-center/"
$ws.Range("C14").Value = 25

# Reflect the final cursor/selection position left by the author
$ws.Range("D15").Select() | Out-Null
